$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I12").Value = "sd"
$ws.Range("J12").Value = "Statement-non-opinion"
$ws.Range("I13").Value = "b"
$ws.Range("J13").Value = "Acknowledge (Backchannel)"
$ws.Range("I30").Value = "aa"
$ws.Range("J30").Value = "Agree/Accept"
$ws.Range("I35").Value = "sd"
$ws.Range("J35").Value = "Statement-non-opinion"
$ws.Range("I42").Value = "sd"
$ws.Range("J42").Value = "Statement-non-opinion"
$ws.Range("I51").Value = "sd"
$ws.Range("J51").Value = "Statement-non-opinion"
$ws.Range("I54").Value = "sd"
$ws.Range("J54").Value = "Statement-non-opinion"
$ws.Range("I64").Value = "sd"
$ws.Range("J64").Value = "Statement-non-opinion"
$ws.Range("I74").Value = "sd"
$ws.Range("J74").Value = "Statement-non-opinion"
$ws.Range("I90").Value = "sd"
$ws.Range("J90").Value = "Statement-non-opinion"
$ws.Range("I98").Value = "sd"
$ws.Range("J98").Value = "Statement-non-opinion"
$ws.Range("I119").Value = "aa"
$ws.Range("J119").Value = "Agree/Accept"
$ws.Range("I120").Value = "aa"
$ws.Range("J120").Value = "Agree/Accept"
$ws.Range("I130").Value = "sv"
$ws.Range("J130").Value = "Statement-opinion"
$ws.Range("I139").Value = "sv"
$ws.Range("J139").Value = "Statement-opinion"
$ws.Range("I149").Value = "sv"
$ws.Range("J149").Value = "Statement-opinion"
$ws.Range("I158").Value = "sd"
$ws.Range("J158").Value = "Statement-non-opinion"
$ws.Range("I181").Value = "aa"
$ws.Range("J181").Value = "Agree/Accept"
$ws.Range("I209").Value = "sv"
$ws.Range("J209").Value = "Statement-opinion"
$ws.Range("I210").Value = "sd"
$ws.Range("J210").Value = "Statement-non-opinion"
$ws.Range("I214").Value = "sv"
$ws.Range("J214").Value = "Statement-opinion"
$ws.Range("I223").Value = "sd"
$ws.Range("J223").Value = "Statement-non-opinion"
$ws.Range("I224").Value = "sd"
$ws.Range("J224").Value = "Statement-non-opinion"
$ws.Range("I233").Value = "ba"
$ws.Range("J233").Value = "Appreciation"
$ws.Range("I234").Value = "b"
$ws.Range("J234").Value = "Acknowledge (Backchannel)"
$ws.Range("I254").Value = "sd"
$ws.Range("J254").Value = "Statement-non-opinion"
$ws.Range("I258").Value = "sv"
$ws.Range("J258").Value = "Statement-opinion"
$ws.Range("I262").Value = "sv"
$ws.Range("J262").Value = "Statement-opinion"
$ws.Range("I269").Value = "aa"
$ws.Range("J269").Value = "Agree/Accept"
$ws.Range("I270").Value = "aa"
$ws.Range("J270").Value = "Agree/Accept"
$ws.Range("I288").Value = "sd"
$ws.Range("J288").Value = "Statement-non-opinion"
$ws.Range("I289").Value = "sd"
$ws.Range("J289").Value = "Statement-non-opinion"
$ws.Range("I295").Value = "sv"
$ws.Range("J295").Value = "Statement-opinion"
$ws.Range("I322").Value = "ba"
$ws.Range("J322").Value = "Appreciation"
$ws.Range("I324").Value = "b"
$ws.Range("J324").Value = "Acknowledge (Backchannel)"
$ws.Range("I325").Value = "sv"
$ws.Range("J325").Value = "Statement-opinion"
$ws.Range("I327").Value = "sd"
$ws.Range("J327").Value = "Statement-non-opinion"
$ws.Range("I331").Value = "sd"
$ws.Range("J331").Value = "Statement-non-opinion"
$ws.Range("I334").Value = "sv"
$ws.Range("J334").Value = "Statement-opinion"
$ws.Range("I348").Value = "aa"
$ws.Range("J348").Value = "Agree/Accept"
$ws.Range("I351").Value = "sv"
$ws.Range("J351").Value = "Statement-opinion"
$ws.Range("I373").Value = "sv"
$ws.Range("J373").Value = "Statement-opinion"
$ws.Range("I379").Value = "sd"
$ws.Range("J379").Value = "Statement-non-opinion"
$ws.Range("I381").Value = "qy"
$ws.Range("J381").Value = "Yes-No-Question"
$ws.Range("I389").Value = "ba"
$ws.Range("J389").Value = "Appreciation"
$ws.Range("I390").Value = "sv"
$ws.Range("J390").Value = "Statement-opinion"
$ws.Range("I394").Value = "ba"
$ws.Range("J394").Value = "Appreciation"
$ws.Range("I401").Value = "sv"
$ws.Range("J401").Value = "Statement-opinion"
$ws.Range("I404").Value = "ba"
$ws.Range("J404").Value = "Appreciation"
$ws.Range("I412").Value = "ba"
$ws.Range("J412").Value = "Appreciation"
$ws.Range("I415").Value = "sd"
$ws.Range("J415").Value = "Statement-non-opinion"
$ws.Range("I421").Value = "sd"
$ws.Range("J421").Value = "Statement-non-opinion"
$ws.Range("I436").Value = "ba"
$ws.Range("J436").Value = "Appreciation"
$ws.Range("I438").Value = "sv"
$ws.Range("J438").Value = "Statement-opinion"
$ws.Range("I451").Value = "sv"
$ws.Range("J451").Value = "Statement-opinion"
$ws.Range("I453").Value = "aa"
$ws.Range("J453").Value = "Agree/Accept"
$ws.Range("I460").Value = "sv"
$ws.Range("J460").Value = "Statement-opinion"
$ws.Range("I475").Value = "aa"
$ws.Range("J475").Value = "Agree/Accept"
$ws.Range("I476").Value = "sv"
$ws.Range("J476").Value = "Statement-opinion"
$ws.Range("I480").Value = "aa"
$ws.Range("J480").Value = "Agree/Accept"
$ws.Range("I507").Value = "sv"
$ws.Range("J507").Value = "Statement-opinion"
$ws.Range("I520").Value = "sd"
$ws.Range("J520").Value = "Statement-non-opinion"
$ws.Range("I525").Value = "aa"
$ws.Range("J525").Value = "Agree/Accept"
$ws.Range("I530").Value = "ba"
$ws.Range("J530").Value = "Appreciation"
$ws.Range("I544").Value = "sd"
$ws.Range("J544").Value = "Statement-non-opinion"
$ws.Range("I549").Value = "sd"
$ws.Range("J549").Value = "Statement-non-opinion"
$ws.Range("I562").Value = "sd"
$ws.Range("J562").Value = "Statement-non-opinion"
$ws.Range("I575").Value = "sd"
$ws.Range("J575").Value = "Statement-non-opinion"
$ws.Range("I583").Value = "sv"
$ws.Range("J583").Value = "Statement-opinion"
